$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.164.79'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').Value = '3.131.78'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.23'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -2.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.60'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -4.93%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '3.123.18'
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('E9').Value = '  -2.56%  '
$ws.Range('E10').Value = '  -3.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.24'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.457'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -2.83%  '
$ws.Range('E13').Value = '  -2.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.22'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -3.23%  '
$ws.Range('D15').Value = '3.644.38'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('E16').Value = '  +3.15%  '
$ws.Range('D17').Value = '63.131.68'
$ws.Range('E17').Value = '  -1.94%  '
$ws.Range('D18').Value = '3.129.42'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.70'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -2.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '473.31'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.23'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -4.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.698'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -2.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.72'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -0.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '86.65'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.00'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -3.77%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.13'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -2.82%  '
$ws.Range('E29').Value = '  -6.04%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.79'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +0.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.109'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -5.79%  '
$ws.Range('E34').Value = '  -4.09%  '
$ws.Range('E35').Value = '  -2.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.82'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -2.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.13'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('D38').Value = '0.0₃0711'
$ws.Range('E38').Value = '  -3.74%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '423.33'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -5.57%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0387'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.24'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('E42').Value = '  -9.81%  '
$ws.Range('D43').Value = '2.891.66'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.113'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -4.17%  '
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('E46').Value = '  -4.30%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.67'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -2.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.28'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -5.45%  '
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.56'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -0.12%  '
